$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 8, shifting rows 8..42 down to 9..43
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new weekly record
$ws.Cells.Item(8, 1).Value = 4
$ws.Cells.Item(8, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(8, 3).Value = "Los Lagos"
$ws.Cells.Item(8, 4).Value = "12/07/2021"
$ws.Cells.Item(8, 5).Value = 10
$ws.Cells.Item(8, 6).Value = 300000000
$ws.Cells.Item(8, 7).Value = "Espárragos"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 480
$ws.Cells.Item(8, 11).Value = 1600
$ws.Cells.Item(8, 12).Value = 1600
$ws.Cells.Item(8, 13).Value = 1600
$ws.Cells.Item(8, 14).Value = "$/kilo"
$ws.Cells.Item(8, 15).Value = "Provincia de Linares"
$ws.Cells.Item(8, 16).Value = 1600
$ws.Cells.Item(8, 17).Value = 1
$ws.Cells.Item(8, 18).Value = "Hortaliza"
